$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("_set_TECHNOLOGIES")
$ws2 = $wb.Worksheets.Item("_set_PRODUCTS")
$ws3 = $wb.Worksheets.Item("_set_AIR_EMISSIONS")

# Drop the "Technology_Assumption" column (column B) from both the
# technologies and products sets - the remaining "Mix" column shifts left
# into column B.
$ws1.Columns("B").Delete()
$ws2.Columns("B").Delete()

# Update the saved selections on every sheet.
$ws2.Range("A2:B8").Select()
$ws3.Range("A2").Select()
$ws1.Range("D14").Select()

# Reset the zoom that had been applied to the technologies sheet.
$excel.ActiveWindow.Zoom = 100

# Make _set_TECHNOLOGIES the active (selected) tab/sheet of the workbook.
$ws1.Activate()
